# Apply updated crypto price/volume values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.214.42'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '1.658.96'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").Value = "'" + '219.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").Value = "'" + '0.5271'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").Value = "'" + '0.2689'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.30%  '
$ws.Range("D9").Value = "'" + '0.06385'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("E10").Value = '  -1.50%  '
$ws.Range("D11").Value = "'" + '0.07688'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.12%  '
$ws.Range("D12").Value = "'" + '4.623'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.10%  '
$ws.Range("D13").Value = '1.664.77'
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").Value = '1.888.03'
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("D15").Value = "'" + '0.5646'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = '0.0₅8273'
$ws.Range("E16").Value = '  +2.10%  '
$ws.Range("D17").Value = "'" + '65.83'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").Value = '26.201.69'
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("D20").Value = "'" + '4.689'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("E21").Value = '  +1.42%  '
$ws.Range("D22").Value = "'" + '191.81'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.34%  '
$ws.Range("D23").Value = "'" + '6.003'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.05%  '
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("D25").Value = "'" + '145.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("D27").Value = "'" + '7.303'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.92%  '
$ws.Range("D28").Value = "'" + '16.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.16%  '
$ws.Range("D29").Value = "'" + '1.527'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").Value = "'" + '0.05655'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.06%  '
$ws.Range("D31").Value = "'" + '1.280'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("D32").Value = "'" + '3.499'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.42%  '
$ws.Range("D33").Value = "'" + '3.401'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.79%  '
$ws.Range("E34").Value = '  -1.08%  '
$ws.Range("D35").Value = "'" + '0.9536'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.33%  '
$ws.Range("D36").Value = "'" + '2.794'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("E37").Value = '  -1.01%  '
$ws.Range("D38").Value = "'" + '0.5776'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.64%  '
$ws.Range("D39").Value = "'" + '0.01605'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("D40").Value = "'" + '5.985'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").Value = "'" + '0.8362'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.55%  '
$ws.Range("D43").Value = '1.030.50'
$ws.Range("E43").Value = '  -4.71%  '
$ws.Range("D44").Value = "'" + '101.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("D45").Value = '1.798.25'
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").Value = "'" + '58.57'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").Value = '0.0₈107'
$ws.Range("E47").Value = '  +4.08%  '
$ws.Range("D48").Value = "'" + '1.005'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("D49").Value = "'" + '0.05346'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.97%  '
$ws.Range("D50").Value = "'" + '8.089'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("E51").Value = '  -1.50%  '
